$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26, pushing existing rows 26-37 down to 27-38.
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the new weekly observation.
$ws.Range("A26").Value = 9
$ws.Range("B26").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C26").Value = 'Metropolitana'
$ws.Range("D26").Value = 44510
$ws.Range("E26").Value = 13
$ws.Range("F26").Value = 100112029
$ws.Range("G26").Value = 'Orégano'
$ws.Range("H26").Value = 'Sin especificar'
$ws.Range("I26").Value = 'Primera'
$ws.Range("J26").Value = 16
$ws.Range("K26").Value = 9000
$ws.Range("L26").Value = 10000
$ws.Range("M26").Value = 9500
$ws.Range("N26").Value = '$/docena de atados'
$ws.Range("O26").Value = 'Región Metropolitana'
$ws.Range("P26").Value = 3167
$ws.Range("Q26").Value = 3
$ws.Range("R26").Value = 'Hortaliza'
